$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "73×86=6278" "16×23=368"
Replace-Text "49×97=4753" "85×85=7225"
Replace-Text "84×92=7728" "41×62=2542"
Replace-Text "76×36=2736" "95×51=4845"
Replace-Text "19×35=665" "34×84=2856"
Replace-Text "42×47=1974" "51×55=2805"
Replace-Text "74×53=3922" "40×51=2040"
Replace-Text "71×63=4473" "48×17=816"
Replace-Text "99×95=9405" "97×85=8245"
Replace-Text "53×30=1590" "19×68=1292"
Replace-Text "51×93=4743" "73×66=4818"
Replace-Text "40×12=480" "41×97=3977"
Replace-Text "94×45=4230" "20×94=1880"
Replace-Text "67×85=5695" "93×60=5580"
Replace-Text "39×66=2574" "53×52=2756"
Replace-Text "46×42=1932" "41×20=820"
Replace-Text "53×32=1696" "46×18=828"
Replace-Text "11×91=1001" "74×78=5772"
Replace-Text "43×25=1075" "49×68=3332"
Replace-Text "36×63=2268" "74×78=5772"
Replace-Text "56×76=4256" "90×23=2070"
Replace-Text "40×46=1840" "52×57=2964"
Replace-Text "72×91=6552" "84×65=5460"
Replace-Text "72×49=3528" "99×85=8415"
Replace-Text "11×97=1067" "11×76=836"

Write-Output "Done applying replacements"
